# feat: add 2022-Q1 data
#
# The workbook has sheets: 2020-Q4, 2021-Q1, 2021-Q2, 2021-Q3, 2021-Q4, 总计
# We need to insert a new "2022-Q1" sheet (fund holdings detail) right before
# the "总计" (totals) sheet, and add a corresponding row to "总计".
#
# To reproduce the exact sheetId / r:id assignment seen in the target
# workbook (2022-Q1 reuses the old "总计" sheetId=6, and a freshly-created
# sheet becomes the new "总计" with sheetId=7), we rename the existing
# "总计" sheet to "2022-Q1" and populate it with the new fund data, then
# add a brand-new sheet, name it "总计", move it to the end, and populate
# it with the updated totals table.

$wb = $excel.ActiveWorkbook

# --- Step 1: rename the existing "总计" sheet to "2022-Q1" -----------------
$oldTotal = $wb.Worksheets.Item("总计")
$oldTotal.Name = "2022-Q1"

# --- Step 2: create the new "总计" sheet and move it to the end ------------
$newTotal = $wb.Worksheets.Add()
$newTotal.Name = "总计"

# Re-fetch sheet references by name right before using them: object
# references captured before sheets are added/reordered can go stale.
$q1 = $wb.Worksheets.Item("2022-Q1")
$newTotal = $wb.Worksheets.Item("总计")
$newTotal.Move($null, $q1)

Write-Output "Sheet order:"
foreach ($s in $wb.Worksheets) {
    Write-Output $s.Name
}

# A sheet that still has the "header row bold+border" / "index column
# bold+border" style (s="2" in the OOXML) that we want to replicate on the
# new/rewritten sheets, so we can copy it across with PasteSpecial.
$styleSource = $wb.Worksheets.Item("2021-Q4")

# =============================================================================
# Step 3: populate "2022-Q1" (fund holdings detail) with the new fund data,
# replacing the old "总计" table that used to live on this sheet/tab.
# =============================================================================
$q1 = $wb.Worksheets.Item("2022-Q1")
$q1.Cells.Clear()

# Copy header-row (bold + border) style from the source sheet's B1:H1 onto
# the new header row, and the index-column style from A2:A6 onto column A.
$styleSource.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)
$styleSource.Range("A2:A6").Copy()
$q1.Range("A2:A6").PasteSpecial(-4122)

# Header row
$q1.Cells.Item(1,2).Value = "基金代码"
$q1.Cells.Item(1,3).Value = "基金名称"
$q1.Cells.Item(1,4).Value = "基金规模"
$q1.Cells.Item(1,5).Value = "股票总仓位"
$q1.Cells.Item(1,6).Value = "仓位占比"
$q1.Cells.Item(1,7).Value = "持有市值(亿元)"
$q1.Cells.Item(1,8).Value = "仓位排名"

# Fund rows: column A is a numeric 0-based row index; B (fund code) and D-G
# (scale/position/value figures) are stored as literal text in the source
# workbook (leading zeros / fixed decimal formatting must be preserved), so
# force text entry via NumberFormat "@" before assigning the string value.
# C (fund name) and H (rank) are plain text / number respectively and don't
# need any special handling.
$q1FundRows = @(
    @(0, "011738", "华安兴安优选一年持有期混合型证券投资基金A", "25.77", "54.03", "3.33", "0.8581", 1),
    @(1, "011739", "华安兴安优选一年持有期混合型证券投资基金C", "10.03", "54.03", "3.33", "0.3340", 1),
    @(2, "011390", "华安添祥6个月持有期混合型证券投资基金",     "8.25",  "33.54", "3.92", "0.3234", 1),
    @(3, "005695", "华安睿明两年定期开放灵活配置混合A",         "1.98",  "93.49", "9.39", "0.1859", 1),
    @(4, "005696", "华安睿明两年定期开放灵活配置混合C",         "0.10",  "93.49", "9.39", "0.0094", 1)
)

$r = 2
foreach ($row in $q1FundRows) {
    $q1.Cells.Item($r, 1).Value = $row[0]

    $q1.Cells.Item($r, 2).NumberFormat = "@"
    $q1.Cells.Item($r, 2).Value = $row[1]

    $q1.Cells.Item($r, 3).Value = $row[2]

    $q1.Cells.Item($r, 4).NumberFormat = "@"
    $q1.Cells.Item($r, 4).Value = $row[3]

    $q1.Cells.Item($r, 5).NumberFormat = "@"
    $q1.Cells.Item($r, 5).Value = $row[4]

    $q1.Cells.Item($r, 6).NumberFormat = "@"
    $q1.Cells.Item($r, 6).Value = $row[5]

    $q1.Cells.Item($r, 7).NumberFormat = "@"
    $q1.Cells.Item($r, 7).Value = $row[6]

    $q1.Cells.Item($r, 8).Value = $row[7]

    $r = $r + 1
}

# =============================================================================
# Step 4: populate the new "总计" sheet with the updated totals table
# (adds the 2022-Q1 row on top, shifting the rest down by one).
# =============================================================================
$total = $wb.Worksheets.Item("总计")
$total.Cells.Clear()

$styleSource.Range("B1:D1").Copy()
$total.Range("B1:D1").PasteSpecial(-4122)
$styleSource.Range("A2:A7").Copy()
$total.Range("A2:A7").PasteSpecial(-4122)

$total.Cells.Item(1,2).Value = "日期"
$total.Cells.Item(1,3).Value = "持有数量(只)"
$total.Cells.Item(1,4).Value = "持有市值(亿元)"

$totalRows = @(
    @(0, "2022-Q1", 5,  1.71),
    @(1, "2021-Q4", 7,  2.06),
    @(2, "2021-Q3", 8,  1.58),
    @(3, "2021-Q2", 9,  0.94),
    @(4, "2021-Q1", 12, 0.94),
    @(5, "2020-Q4", 2,  0.09)
)

$r = 2
foreach ($row in $totalRows) {
    $total.Cells.Item($r, 1).Value = $row[0]
    $total.Cells.Item($r, 2).Value = $row[1]
    $total.Cells.Item($r, 3).Value = $row[2]
    $total.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}
